$wb = $excel.ActiveWorkbook

# ----- Overview sheet: update cell values -----
$ws = $wb.Worksheets.Item('Overview')
$ws.Range('A2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.md'
$ws.Range('B2').Value = 'e2e\883edf0d-92c2-4420-9325-3832dfe69631.md'
$ws.Range('G2').Value = '2016-09-06 04:30:28'
$ws.Range('A3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.md'
$ws.Range('B3').Value = 'e2e\a2862e86-1707-4973-b8a8-c4401aefcc61.md'
$ws.Range('G3').Value = '2016-09-06 04:30:28'
$ws.Range('A4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.md'
$ws.Range('B4').Value = 'e2e\320b5078-4a44-4936-8350-bf022276e5ff.md'
$ws.Range('E4').Value = 'Ready for handoff'
$ws.Range('F4').Value = 'Ready for handoff'
$ws.Range('G4').Value = '2016-09-06 04:33:54'
$ws.Range('A5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.md'
$ws.Range('B5').Value = 'e2e\551db600-cdcf-4134-9279-bc5d22b2881d.md'
$ws.Range('E5').Value = 'Ready for handoff'
$ws.Range('F5').Value = 'Ready for handoff'
$ws.Range('G5').Value = '2016-09-06 04:33:54'

# ----- Overview sheet: rebuild hyperlinks in new order -----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('B2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md', "", "", 'e2e\883edf0d-92c2-4420-9325-3832dfe69631.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('B3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md', "", "", 'e2e\a2862e86-1707-4973-b8a8-c4401aefcc61.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('B4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ace7d096139661ed1bc8fd4d36cd0c744ef42239/e2e/883edf0d-92c2-4420-9325-3832dfe69631.md', "", "", 'e2e\320b5078-4a44-4936-8350-bf022276e5ff.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('B5'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ace7d096139661ed1bc8fd4d36cd0c744ef42239/e2e/a2862e86-1707-4973-b8a8-c4401aefcc61.md', "", "", 'e2e\551db600-cdcf-4134-9279-bc5d22b2881d.md') | Out-Null

# ----- zh-cn sheet: update cell values -----
$ws = $wb.Worksheets.Item('zh-cn')
$ws.Range('A2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.md'
$ws.Range('G2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.16bc7cd6488f4306c8d424624b1972ef07a53310.zh-cn.xlf'
$ws.Range('H2').Value = '2016-09-06 04:30:22'
$ws.Range('I2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.md'
$ws.Range('J2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.16bc7cd6488f4306c8d424624b1972ef07a53310.zh-cn.xlf'
$ws.Range('K2').Value = '2016-09-06 04:31:43'
$ws.Range('A3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.md'
$ws.Range('G3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.ea431df66d15f0d6b57cdd90889c8d996b240be8.zh-cn.xlf'
$ws.Range('H3').Value = '2016-09-06 04:30:22'
$ws.Range('I3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.md'
$ws.Range('J3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.ea431df66d15f0d6b57cdd90889c8d996b240be8.zh-cn.xlf'
$ws.Range('K3').Value = '2016-09-06 04:31:43'
$ws.Range('A4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.md'
$ws.Range('C4').Value = 'Ready for handoff'
$ws.Range('G4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.b6ac3c0588fe2f63ab6785bd3be6eaf6a27eae14.zh-cn.xlf'
$ws.Range('H4').Value = '2016-09-06 04:33:48'
$ws.Range('I4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.md'
$ws.Range('J4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.b6ac3c0588fe2f63ab6785bd3be6eaf6a27eae14.zh-cn.xlf'
$ws.Range('K4').Value = '2016-09-06 04:33:02'
$ws.Range('P4').Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41bf586e87b3c93d3885f6c52bbe3dd1f64c39c2/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md.'
$ws.Range('A5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.md'
$ws.Range('C5').Value = 'Ready for handoff'
$ws.Range('G5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.db1ef45905a064d61c1c1fd588a9a6181b0ac0f6.zh-cn.xlf'
$ws.Range('H5').Value = '2016-09-06 04:33:48'
$ws.Range('I5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.md'
$ws.Range('J5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.db1ef45905a064d61c1c1fd588a9a6181b0ac0f6.zh-cn.xlf'
$ws.Range('K5').Value = '2016-09-06 04:33:02'
$ws.Range('P5').Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41bf586e87b3c93d3885f6c52bbe3dd1f64c39c2/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md.'

# ----- zh-cn sheet: rebuild hyperlinks in new order -----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md', "", "", '883edf0d-92c2-4420-9325-3832dfe69631.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e298fd23f0fd68a380dc8646364d036106702886/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md', "", "", '883edf0d-92c2-4420-9325-3832dfe69631.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md', "", "", 'a2862e86-1707-4973-b8a8-c4401aefcc61.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e298fd23f0fd68a380dc8646364d036106702886/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md', "", "", 'a2862e86-1707-4973-b8a8-c4401aefcc61.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ace7d096139661ed1bc8fd4d36cd0c744ef42239/e2e/883edf0d-92c2-4420-9325-3832dfe69631.md', "", "", '320b5078-4a44-4936-8350-bf022276e5ff.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7f71d378323142611d8f28227b969d3e0eff36e/e2e/883edf0d-92c2-4420-9325-3832dfe69631.md', "", "", '320b5078-4a44-4936-8350-bf022276e5ff.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ace7d096139661ed1bc8fd4d36cd0c744ef42239/e2e/a2862e86-1707-4973-b8a8-c4401aefcc61.md', "", "", '551db600-cdcf-4134-9279-bc5d22b2881d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I5'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7f71d378323142611d8f28227b969d3e0eff36e/e2e/a2862e86-1707-4973-b8a8-c4401aefcc61.md', "", "", '551db600-cdcf-4134-9279-bc5d22b2881d.md') | Out-Null

# ----- de-de sheet: update cell values -----
$ws = $wb.Worksheets.Item('de-de')
$ws.Range('A2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.md'
$ws.Range('G2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.16bc7cd6488f4306c8d424624b1972ef07a53310.de-de.xlf'
$ws.Range('H2').Value = '2016-09-06 04:30:28'
$ws.Range('I2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.md'
$ws.Range('J2').Value = '883edf0d-92c2-4420-9325-3832dfe69631.16bc7cd6488f4306c8d424624b1972ef07a53310.de-de.xlf'
$ws.Range('K2').Value = '2016-09-06 04:31:50'
$ws.Range('A3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.md'
$ws.Range('G3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.ea431df66d15f0d6b57cdd90889c8d996b240be8.de-de.xlf'
$ws.Range('H3').Value = '2016-09-06 04:30:28'
$ws.Range('I3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.md'
$ws.Range('J3').Value = 'a2862e86-1707-4973-b8a8-c4401aefcc61.ea431df66d15f0d6b57cdd90889c8d996b240be8.de-de.xlf'
$ws.Range('K3').Value = '2016-09-06 04:31:50'
$ws.Range('A4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.md'
$ws.Range('C4').Value = 'Ready for handoff'
$ws.Range('G4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.b6ac3c0588fe2f63ab6785bd3be6eaf6a27eae14.de-de.xlf'
$ws.Range('H4').Value = '2016-09-06 04:33:54'
$ws.Range('I4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.md'
$ws.Range('J4').Value = '320b5078-4a44-4936-8350-bf022276e5ff.b6ac3c0588fe2f63ab6785bd3be6eaf6a27eae14.de-de.xlf'
$ws.Range('K4').Value = '2016-09-06 04:33:17'
$ws.Range('P4').Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41bf586e87b3c93d3885f6c52bbe3dd1f64c39c2/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md.'
$ws.Range('A5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.md'
$ws.Range('C5').Value = 'Ready for handoff'
$ws.Range('G5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.db1ef45905a064d61c1c1fd588a9a6181b0ac0f6.de-de.xlf'
$ws.Range('H5').Value = '2016-09-06 04:33:54'
$ws.Range('I5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.md'
$ws.Range('J5').Value = '551db600-cdcf-4134-9279-bc5d22b2881d.db1ef45905a064d61c1c1fd588a9a6181b0ac0f6.de-de.xlf'
$ws.Range('K5').Value = '2016-09-06 04:33:17'
$ws.Range('P5').Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41bf586e87b3c93d3885f6c52bbe3dd1f64c39c2/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md.'

# ----- de-de sheet: rebuild hyperlinks in new order -----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md', "", "", '883edf0d-92c2-4420-9325-3832dfe69631.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ca7c51d13ebfdd23f376fa4e73158577a6155ccf/e2e/320b5078-4a44-4936-8350-bf022276e5ff.md', "", "", '883edf0d-92c2-4420-9325-3832dfe69631.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1e09372df43343b57465c257c03a856a5254ae5/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md', "", "", 'a2862e86-1707-4973-b8a8-c4401aefcc61.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ca7c51d13ebfdd23f376fa4e73158577a6155ccf/e2e/551db600-cdcf-4134-9279-bc5d22b2881d.md', "", "", 'a2862e86-1707-4973-b8a8-c4401aefcc61.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ace7d096139661ed1bc8fd4d36cd0c744ef42239/e2e/883edf0d-92c2-4420-9325-3832dfe69631.md', "", "", '320b5078-4a44-4936-8350-bf022276e5ff.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/aa444ec18ec09cc4fdaec889f215d11e1cc5814d/e2e/883edf0d-92c2-4420-9325-3832dfe69631.md', "", "", '320b5078-4a44-4936-8350-bf022276e5ff.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ace7d096139661ed1bc8fd4d36cd0c744ef42239/e2e/a2862e86-1707-4973-b8a8-c4401aefcc61.md', "", "", '551db600-cdcf-4134-9279-bc5d22b2881d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('I5'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/aa444ec18ec09cc4fdaec889f215d11e1cc5814d/e2e/a2862e86-1707-4973-b8a8-c4401aefcc61.md', "", "", '551db600-cdcf-4134-9279-bc5d22b2881d.md') | Out-Null

# ----- Column width adjustments (column P widened to 40) -----
$wsZh = $wb.Worksheets.Item('zh-cn')
$wsZh.Columns.Item(16).ColumnWidth = 40
$wsDe = $wb.Worksheets.Item('de-de')
$wsDe.Columns.Item(16).ColumnWidth = 40
